$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 9-18 (row 19 "Klay Thompson / Dallas Mavericks" removed)
$data = @(
    @("Corey Kispert",   "SG,SF", "Washington Wizards"),
    @("Rudy Gobert",     "C",     "Minnesota Timberwolves"),
    @("Jakob Poeltl",    "C",     "Toronto Raptors"),
    @("Nikola Jokic",    "C",     "Denver Nuggets"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Clint Capela",    "C",     "Atlanta Hawks"),
    @("Jerami Grant",    "SF,PF", "Portland Trail Blazers"),
    @("Jaylen Brown",    "SG,SF", "Boston Celtics"),
    @("Paolo Banchero",  "SF,PF", "Orlando Magic"),
    @("Chet Holmgren",   "PF,C",  "Oklahoma City Thunder")
)

$row = 9
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

# Remove the now-extra last row (previously row 19, "Klay Thompson")
$ws.Range("A19:C19").Delete()
